$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Apoe"
$ws.Cells.Item(2, 3).Value = "Vldlr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 29.32133366666666
$ws.Cells.Item(2, 8).Value = 87.964001
$ws.Cells.Item(2, 9).Value = 0.006401919837078288
$ws.Cells.Item(2, 10).Value = 0.006401919837078288
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.463191
$ws.Cells.Item(2, 14).Value = 1.389573
$ws.Cells.Item(2, 15).Value = 0.0353316468093919
$ws.Cells.Item(2, 16).Value = 0.0353316468093919
$ws.Cells.Item(2, 17).Value = 13.581377862397
$ws.Cells.Item(2, 18).Value = 122.232400761573
$ws.Cells.Item(2, 19).Value = 0.0002261903705856898
$ws.Cells.Item(2, 20).Value = 0.0002261903705856898

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Apoe"
$ws.Cells.Item(3, 3).Value = "Vldlr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 29.32133366666666
$ws.Cells.Item(3, 8).Value = 87.964001
$ws.Cells.Item(3, 9).Value = 0.006401919837078288
$ws.Cells.Item(3, 10).Value = 0.006401919837078288
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.932929333333332
$ws.Cells.Item(3, 14).Value = 29.798788
$ws.Cells.Item(3, 15).Value = 0.7576717833204485
$ws.Cells.Item(3, 16).Value = 0.7576717833204486
$ws.Cells.Item(3, 17).Value = 291.2467352700875
$ws.Cells.Item(3, 18).Value = 2621.220617430788
$ws.Cells.Item(3, 19).Value = 0.004850554019633661
$ws.Cells.Item(3, 20).Value = 0.004850554019633663

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Apoe"
$ws.Cells.Item(4, 3).Value = "Vldlr"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 29.32133366666666
$ws.Cells.Item(4, 8).Value = 87.964001
$ws.Cells.Item(4, 9).Value = 0.006401919837078288
$ws.Cells.Item(4, 10).Value = 0.006401919837078288
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1285543333333333
$ws.Cells.Item(4, 14).Value = 0.385663
$ws.Cells.Item(4, 15).Value = 0.009805968382697785
$ws.Cells.Item(4, 16).Value = 0.009805968382697785
$ws.Cells.Item(4, 17).Value = 3.769384501962555
$ws.Cells.Item(4, 18).Value = 33.92446051766299
$ws.Cells.Item(4, 19).Value = 0.00006277702351095544
$ws.Cells.Item(4, 20).Value = 0.00006277702351095546

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Apoe"
$ws.Cells.Item(5, 3).Value = "Vldlr"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 29.32133366666666
$ws.Cells.Item(5, 8).Value = 87.964001
$ws.Cells.Item(5, 9).Value = 0.006401919837078288
$ws.Cells.Item(5, 10).Value = 0.006401919837078288
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.585130333333333
$ws.Cells.Item(5, 14).Value = 7.755391
$ws.Cells.Item(5, 15).Value = 0.1971906014874617
$ws.Cells.Item(5, 16).Value = 0.1971906014874618
$ws.Cells.Item(5, 17).Value = 75.79946907548788
$ws.Cells.Item(5, 18).Value = 682.195221679391
$ws.Cells.Item(5, 19).Value = 0.00126239842334798
$ws.Cells.Item(5, 20).Value = 0.001262398423347981

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Apoe"
$ws.Cells.Item(6, 3).Value = "Vldlr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 45.524413
$ws.Cells.Item(6, 8).Value = 136.573239
$ws.Cells.Item(6, 9).Value = 0.009939644832300594
$ws.Cells.Item(6, 10).Value = 0.009939644832300592
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.463191
$ws.Cells.Item(6, 14).Value = 1.389573
$ws.Cells.Item(6, 15).Value = 0.0353316468093919
$ws.Cells.Item(6, 16).Value = 0.0353316468093919
$ws.Cells.Item(6, 17).Value = 21.086498381883
$ws.Cells.Item(6, 18).Value = 189.778485436947
$ws.Cells.Item(6, 19).Value = 0.000351184020625642
$ws.Cells.Item(6, 20).Value = 0.0003511840206256419

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Apoe"
$ws.Cells.Item(7, 3).Value = "Vldlr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 45.524413
$ws.Cells.Item(7, 8).Value = 136.573239
$ws.Cells.Item(7, 9).Value = 0.009939644832300594
$ws.Cells.Item(7, 10).Value = 0.009939644832300592
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.932929333333332
$ws.Cells.Item(7, 14).Value = 29.798788
$ws.Cells.Item(7, 15).Value = 0.7576717833204485
$ws.Cells.Item(7, 16).Value = 0.7576717833204486
$ws.Cells.Item(7, 17).Value = 452.1907772704813
$ws.Cells.Item(7, 18).Value = 4069.716995434332
$ws.Cells.Item(7, 19).Value = 0.007530988425661071
$ws.Cells.Item(7, 20).Value = 0.007530988425661071

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Apoe"
$ws.Cells.Item(8, 3).Value = "Vldlr"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 45.524413
$ws.Cells.Item(8, 8).Value = 136.573239
$ws.Cells.Item(8, 9).Value = 0.009939644832300594
$ws.Cells.Item(8, 10).Value = 0.009939644832300592
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1285543333333333
$ws.Cells.Item(8, 14).Value = 0.385663
$ws.Cells.Item(8, 15).Value = 0.009805968382697785
$ws.Cells.Item(8, 16).Value = 0.009805968382697785
$ws.Cells.Item(8, 17).Value = 5.852360563606333
$ws.Cells.Item(8, 18).Value = 52.671245072457
$ws.Cells.Item(8, 19).Value = 0.00009746784296078506
$ws.Cells.Item(8, 20).Value = 0.00009746784296078503

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Apoe"
$ws.Cells.Item(9, 3).Value = "Vldlr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 45.524413
$ws.Cells.Item(9, 8).Value = 136.573239
$ws.Cells.Item(9, 9).Value = 0.009939644832300594
$ws.Cells.Item(9, 10).Value = 0.009939644832300592
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.585130333333333
$ws.Cells.Item(9, 14).Value = 7.755391
$ws.Cells.Item(9, 15).Value = 0.1971906014874617
$ws.Cells.Item(9, 16).Value = 0.1971906014874618
$ws.Cells.Item(9, 17).Value = 117.6865409534943
$ws.Cells.Item(9, 18).Value = 1059.178868581449
$ws.Cells.Item(9, 19).Value = 0.001960004543053095
$ws.Cells.Item(9, 20).Value = 0.001960004543053095

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Apoe"
$ws.Cells.Item(10, 3).Value = "Vldlr"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4438.215250666667
$ws.Cells.Item(10, 8).Value = 13314.645752
$ws.Cells.Item(10, 9).Value = 0.9690247577915309
$ws.Cells.Item(10, 10).Value = 0.9690247577915307
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.463191
$ws.Cells.Item(10, 14).Value = 1.389573
$ws.Cells.Item(10, 15).Value = 0.0353316468093919
$ws.Cells.Item(10, 16).Value = 0.0353316468093919
$ws.Cells.Item(10, 17).Value = 2055.741360171544
$ws.Cells.Item(10, 18).Value = 18501.6722415439
$ws.Cells.Item(10, 19).Value = 0.0342372404918469
$ws.Cells.Item(10, 20).Value = 0.0342372404918469

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Apoe"
$ws.Cells.Item(11, 3).Value = "Vldlr"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4438.215250666667
$ws.Cells.Item(11, 8).Value = 13314.645752
$ws.Cells.Item(11, 9).Value = 0.9690247577915309
$ws.Cells.Item(11, 10).Value = 0.9690247577915307
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 9.932929333333332
$ws.Cells.Item(11, 14).Value = 29.798788
$ws.Cells.Item(11, 15).Value = 0.7576717833204485
$ws.Cells.Item(11, 16).Value = 0.7576717833204486
$ws.Cells.Item(11, 17).Value = 44084.47845099428
$ws.Cells.Item(11, 18).Value = 396760.3060589486
$ws.Cells.Item(11, 19).Value = 0.7342027163175748
$ws.Cells.Item(11, 20).Value = 0.7342027163175748

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Apoe"
$ws.Cells.Item(12, 3).Value = "Vldlr"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4438.215250666667
$ws.Cells.Item(12, 8).Value = 13314.645752
$ws.Cells.Item(12, 9).Value = 0.9690247577915309
$ws.Cells.Item(12, 10).Value = 0.9690247577915307
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.1285543333333333
$ws.Cells.Item(12, 14).Value = 0.385663
$ws.Cells.Item(12, 15).Value = 0.009805968382697785
$ws.Cells.Item(12, 16).Value = 0.009805968382697785
$ws.Cells.Item(12, 17).Value = 570.5518027392862
$ws.Cells.Item(12, 18).Value = 5134.966224653575
$ws.Cells.Item(12, 19).Value = 0.00950222613695513
$ws.Cells.Item(12, 20).Value = 0.00950222613695513

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Apoe"
$ws.Cells.Item(13, 3).Value = "Vldlr"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4438.215250666667
$ws.Cells.Item(13, 8).Value = 13314.645752
$ws.Cells.Item(13, 9).Value = 0.9690247577915309
$ws.Cells.Item(13, 10).Value = 0.9690247577915307
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.585130333333333
$ws.Cells.Item(13, 14).Value = 7.755391
$ws.Cells.Item(13, 15).Value = 0.1971906014874617
$ws.Cells.Item(13, 16).Value = 0.1971906014874618
$ws.Cells.Item(13, 17).Value = 11473.36487036101
$ws.Cells.Item(13, 18).Value = 103260.283833249
$ws.Cells.Item(13, 19).Value = 0.1910825748451539
$ws.Cells.Item(13, 20).Value = 0.1910825748451539

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Apoe"
$ws.Cells.Item(14, 3).Value = "Vldlr"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 67.02347933333333
$ws.Cells.Item(14, 8).Value = 201.070438
$ws.Cells.Item(14, 9).Value = 0.01463367753909034
$ws.Cells.Item(14, 10).Value = 0.01463367753909034
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.463191
$ws.Cells.Item(14, 14).Value = 1.389573
$ws.Cells.Item(14, 15).Value = 0.0353316468093919
$ws.Cells.Item(14, 16).Value = 0.0353316468093919
$ws.Cells.Item(14, 17).Value = 31.04467241588599
$ws.Cells.Item(14, 18).Value = 279.402051742974
$ws.Cells.Item(14, 19).Value = 0.0005170319263336711
$ws.Cells.Item(14, 20).Value = 0.0005170319263336712

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Apoe"
$ws.Cells.Item(15, 3).Value = "Vldlr"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 67.02347933333333
$ws.Cells.Item(15, 8).Value = 201.070438
$ws.Cells.Item(15, 9).Value = 0.01463367753909034
$ws.Cells.Item(15, 10).Value = 0.01463367753909034
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 9.932929333333332
$ws.Cells.Item(15, 14).Value = 29.798788
$ws.Cells.Item(15, 15).Value = 0.7576717833204485
$ws.Cells.Item(15, 16).Value = 0.7576717833204486
$ws.Cells.Item(15, 17).Value = 665.7394838921269
$ws.Cells.Item(15, 18).Value = 5991.655355029144
$ws.Cells.Item(15, 19).Value = 0.01108752455757897
$ws.Cells.Item(15, 20).Value = 0.01108752455757897

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Apoe"
$ws.Cells.Item(16, 3).Value = "Vldlr"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 67.02347933333333
$ws.Cells.Item(16, 8).Value = 201.070438
$ws.Cells.Item(16, 9).Value = 0.01463367753909034
$ws.Cells.Item(16, 10).Value = 0.01463367753909034
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1285543333333333
$ws.Cells.Item(16, 14).Value = 0.385663
$ws.Cells.Item(16, 15).Value = 0.009805968382697785
$ws.Cells.Item(16, 16).Value = 0.009805968382697785
$ws.Cells.Item(16, 17).Value = 8.61615870337711
$ws.Cells.Item(16, 18).Value = 77.54542833039399
$ws.Cells.Item(16, 19).Value = 0.0001434973792709146
$ws.Cells.Item(16, 20).Value = 0.0001434973792709146

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Apoe"
$ws.Cells.Item(17, 3).Value = "Vldlr"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 67.02347933333333
$ws.Cells.Item(17, 8).Value = 201.070438
$ws.Cells.Item(17, 9).Value = 0.01463367753909034
$ws.Cells.Item(17, 10).Value = 0.01463367753909034
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 2.585130333333333
$ws.Cells.Item(17, 14).Value = 7.755391
$ws.Cells.Item(17, 15).Value = 0.1971906014874617
$ws.Cells.Item(17, 16).Value = 0.1971906014874618
$ws.Cells.Item(17, 17).Value = 173.2644294701398
$ws.Cells.Item(17, 18).Value = 1559.379865231258
$ws.Cells.Item(17, 19).Value = 0.002885623675906782
$ws.Cells.Item(17, 20).Value = 0.002885623675906783
